# Update "countries & provincias Spain" data (paises.xlsx)
# - Move Estonia's row to directly follow Mali (row 136), pushing Bahamas and
#   Sudan del Sur down by one row each (Aruba keeps its row), and refresh
#   Estonia's statistics.
# - Refresh the COVID-19 case statistics for a number of other countries.
# - Update the "Datos actualizados a ..." timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp string (A1) -------------------------------------------------
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 8 de Septiembre de 2020 a las 12:48"

# --- Reorder Mali / Estonia / Bahamas / Sudan del Sur / Aruba --------------
# Row 135 (Mali) and row 139 (Aruba) are unchanged.
# Row 136 becomes Estonia (new figures), row 137 becomes Bahamas (old row-136
# figures) and row 138 becomes Sudan del Sur (old row-137 figures).

# Row 136 -> Estonia (refreshed figures)
$ws.Cells.Item(136, 1).Value = "Estonia"
$ws.Cells.Item(136, 2).Value = 2564
$ws.Cells.Item(136, 3).Value = 32
$ws.Cells.Item(136, 4).Value = 2195
$ws.Cells.Item(136, 5).Value = 305
$ws.Cells.Item(136, 6).Value = 0
$ws.Cells.Item(136, 7).Value = 0
$ws.Cells.Item(136, 8).Value = 64

# Row 137 -> Bahamas (figures that used to live in row 136)
$ws.Cells.Item(137, 1).Value = "Bahamas"
$ws.Cells.Item(137, 2).Value = 2546
$ws.Cells.Item(137, 3).Value = 0
$ws.Cells.Item(137, 4).Value = 976
$ws.Cells.Item(137, 5).Value = 1512
$ws.Cells.Item(137, 6).Value = 0
$ws.Cells.Item(137, 7).Value = 0
$ws.Cells.Item(137, 8).Value = 58

# Row 138 -> Sudan del Sur (figures that used to live in row 137)
$ws.Cells.Item(138, 1).Value = "Sudan del Sur"
$ws.Cells.Item(138, 2).Value = 2545
$ws.Cells.Item(138, 3).Value = 0
$ws.Cells.Item(138, 4).Value = 1290
$ws.Cells.Item(138, 5).Value = 1207
$ws.Cells.Item(138, 6).Value = 0
$ws.Cells.Item(138, 7).Value = 0
$ws.Cells.Item(138, 8).Value = 48

# --- Refresh case counts for other countries -------------------------------

# Row 4 - Estados Unidos
$ws.Cells.Item(4, 2).Value = 6485708
$ws.Cells.Item(4, 3).Value = 133
$ws.Cells.Item(4, 5).Value = 2533543
$ws.Cells.Item(4, 7).Value = 2
$ws.Cells.Item(4, 8).Value = 193536

# Row 7 - Rusia
$ws.Cells.Item(7, 2).Value = 1035789
$ws.Cells.Item(7, 3).Value = 5099
$ws.Cells.Item(7, 4).Value = 850049
$ws.Cells.Item(7, 5).Value = 167747
$ws.Cells.Item(7, 7).Value = 122
$ws.Cells.Item(7, 8).Value = 17993

# Row 38 - Kuwait
$ws.Cells.Item(38, 2).Value = 91244
$ws.Cells.Item(38, 3).Value = 857
$ws.Cells.Item(38, 4).Value = 81654
$ws.Cells.Item(38, 5).Value = 9042
$ws.Cells.Item(38, 7).Value = 2
$ws.Cells.Item(38, 8).Value = 548

# Row 49 - Polonia
$ws.Cells.Item(49, 2).Value = 71526
$ws.Cells.Item(49, 3).Value = 400
$ws.Cells.Item(49, 5).Value = 13480
$ws.Cells.Item(49, 7).Value = 12
$ws.Cells.Item(49, 8).Value = 2136

# Row 61 - Armenia
$ws.Cells.Item(61, 2).Value = 44953
$ws.Cells.Item(61, 3).Value = 108
$ws.Cells.Item(61, 4).Value = 40592
$ws.Cells.Item(61, 5).Value = 3458
$ws.Cells.Item(61, 7).Value = 3
$ws.Cells.Item(61, 8).Value = 903

# Row 66 - Afganistan
$ws.Cells.Item(66, 2).Value = 38520
$ws.Cells.Item(66, 3).Value = 26
$ws.Cells.Item(66, 4).Value = 30715
$ws.Cells.Item(66, 5).Value = 6387
$ws.Cells.Item(66, 7).Value = 3
$ws.Cells.Item(66, 8).Value = 1418

# Row 75 - Australia
$ws.Cells.Item(75, 2).Value = 26374
$ws.Cells.Item(75, 3).Value = 52
$ws.Cells.Item(75, 4).Value = 22724
$ws.Cells.Item(75, 5).Value = 2880

# Row 96 - Malasia
$ws.Cells.Item(96, 2).Value = 9559
$ws.Cells.Item(96, 3).Value = 100
$ws.Cells.Item(96, 4).Value = 9136
$ws.Cells.Item(96, 5).Value = 295

# Row 104 - Finlandia
$ws.Cells.Item(104, 2).Value = 8337
$ws.Cells.Item(104, 3).Value = 10
$ws.Cells.Item(104, 5).Value = 651

# Row 179 - Islas Feroe
$ws.Cells.Item(179, 2).Value = 414
$ws.Cells.Item(179, 3).Value = 1
$ws.Cells.Item(179, 5).Value = 5
